# Apply updated probability values to the team-specific matrix sheet
# (Sheet1 / Bucknell_A.xlsx equivalent) - reflects new simulated game counts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.16
$ws.Range("C2").Value = 0.62
$ws.Range("J2").Value = 0.012
$ws.Range("P2").Value = 0.116
$ws.Range("S2").Value = 0.092
$ws.Range("B3").Value = 0.00641025641025641
$ws.Range("C3").Value = 0.01282051282051282
$ws.Range("J3").Value = 0.03205128205128205
$ws.Range("P3").Value = 0.7243589743589743
$ws.Range("S3").Value = 0.2243589743589744
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("O4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.5666666666666667
$ws.Range("S4").Value = 0.3666666666666666
$ws.Range("B6").Value = 0.0695970695970696
$ws.Range("D6").Value = 0.01465201465201465
$ws.Range("F6").Value = 0.1062271062271062
$ws.Range("J6").Value = 0.1465201465201465
$ws.Range("O6").Value = 0.01465201465201465
$ws.Range("Q6").Value = 0.1721611721611722
$ws.Range("R6").Value = 0.09523809523809523
$ws.Range("S6").Value = 0.3809523809523809
$ws.Range("B7").Value = 0.1270718232044199
$ws.Range("D7").Value = 0.01657458563535912
$ws.Range("F7").Value = 0.09944751381215469
$ws.Range("J7").Value = 0.1325966850828729
$ws.Range("O7").Value = 0.01657458563535912
$ws.Range("Q7").Value = 0.1988950276243094
$ws.Range("R7").Value = 0.08287292817679558
$ws.Range("S7").Value = 0.3259668508287293
$ws.Range("B8").Value = 0.0951219512195122
$ws.Range("D8").Value = 0.007317073170731708
$ws.Range("E8").Value = 0.002439024390243902
$ws.Range("F8").Value = 0.07317073170731707
$ws.Range("J8").Value = 0.1121951219512195
$ws.Range("O8").Value = 0.02439024390243903
$ws.Range("Q8").Value = 0.1365853658536585
$ws.Range("R8").Value = 0.1121951219512195
$ws.Range("S8").Value = 0.4365853658536585
$ws.Range("B9").Value = 0.05045871559633028
$ws.Range("D9").Value = 0.01376146788990826
$ws.Range("F9").Value = 0.07798165137614679
$ws.Range("J9").Value = 0.1284403669724771
$ws.Range("O9").Value = 0.04587155963302753
$ws.Range("Q9").Value = 0.1467889908256881
$ws.Range("R9").Value = 0.1009174311926606
$ws.Range("S9").Value = 0.4357798165137615
$ws.Range("B10").Value = 0.09730668983492616
$ws.Range("D10").Value = 0.01824500434404865
$ws.Range("F10").Value = 0.09904430929626412
$ws.Range("J10").Value = 0.08861859252823631
$ws.Range("O10").Value = 0.02432667245873154
$ws.Range("Q10").Value = 0.1702867072111208
$ws.Range("R10").Value = 0.09904430929626412
$ws.Range("S10").Value = 0.4031277150304083
$ws.Range("G11").Value = 0.1473684210526316
$ws.Range("J11").Value = 0.1087719298245614
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5333333333333333
$ws.Range("S11").Value = 0.01052631578947368
$ws.Range("G12").Value = 0.7468354430379747
$ws.Range("J12").Value = 0.1708860759493671
$ws.Range("K12").Value = 0.006329113924050633
$ws.Range("L12").Value = 0.03164556962025317
$ws.Range("S12").Value = 0.04430379746835443
$ws.Range("G13").Value = 0.5957446808510638
$ws.Range("J13").Value = 0.3617021276595745
$ws.Range("S13").Value = 0.0425531914893617
$ws.Range("F15").Value = 0.007782101167315175
$ws.Range("H15").Value = 0.09727626459143969
$ws.Range("I15").Value = 0.11284046692607
$ws.Range("J15").Value = 0.3385214007782101
$ws.Range("K15").Value = 0.04669260700389105
$ws.Range("M15").Value = 0.01945525291828794
$ws.Range("O15").Value = 0.1011673151750973
$ws.Range("S15").Value = 0.2762645914396887
$ws.Range("F16").Value = 0.03246753246753246
$ws.Range("H16").Value = 0.1688311688311688
$ws.Range("I16").Value = 0.08441558441558442
$ws.Range("J16").Value = 0.3831168831168831
$ws.Range("K16").Value = 0.08441558441558442
$ws.Range("M16").Value = 0.05194805194805195
$ws.Range("O16").Value = 0.05844155844155844
$ws.Range("S16").Value = 0.1363636363636364
$ws.Range("F17").Value = 0.02472527472527472
$ws.Range("H17").Value = 0.1620879120879121
$ws.Range("I17").Value = 0.08791208791208792
$ws.Range("J17").Value = 0.3983516483516483
$ws.Range("K17").Value = 0.1153846153846154
$ws.Range("M17").Value = 0.02197802197802198
$ws.Range("O17").Value = 0.06043956043956044
$ws.Range("S17").Value = 0.1291208791208791
$ws.Range("F18").Value = 0.009049773755656109
$ws.Range("H18").Value = 0.1719457013574661
$ws.Range("I18").Value = 0.1085972850678733
$ws.Range("J18").Value = 0.3891402714932127
$ws.Range("K18").Value = 0.09954751131221719
$ws.Range("M18").Value = 0.01809954751131222
$ws.Range("N18").Value = 0.004524886877828055
$ws.Range("O18").Value = 0.06787330316742081
$ws.Range("S18").Value = 0.1312217194570136
$ws.Range("F19").Value = 0.01341752170481452
$ws.Range("H19").Value = 0.2052091554853986
$ws.Range("I19").Value = 0.09550118389897395
$ws.Range("J19").Value = 0.372533543804262
$ws.Range("K19").Value = 0.1081294396211523
$ws.Range("M19").Value = 0.01736385161799527
$ws.Range("O19").Value = 0.08366219415943173
$ws.Range("S19").Value = 0.1026045777426993
